$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New "K" (strikeout) column values, replacing the old "Strike#" values,
# for rows 2-38 in column G.
$kValues = @{
    2  = 0
    3  = 1
    4  = 1
    5  = 2
    6  = 1
    7  = 0
    8  = 0
    9  = 3
    10 = 0
    11 = 3
    12 = 1
    13 = 1
    14 = 2
    15 = 3
    16 = 2
    17 = 2
    18 = 1
    19 = 1
    20 = 3
    21 = 2
    22 = 0
    23 = 0
    24 = 0
    25 = 0
    26 = 0
    27 = 0
    28 = 0
    29 = 2
    30 = 2
    31 = 2
    32 = 1
    33 = 1
    34 = 0
    35 = 2
    36 = 1
    37 = 1
    38 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
